# Auto-generated edit script: updates cached numeric values in the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR profit-tracking sheets
# per the scheduled price-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1181.125
$ws.Range("I32").Value = 887.5
$ws.Range("J32").Value = 1474.75
$ws.Range("K32").Value = 887.5
$ws.Range("L32").Value = 1474.75
$ws.Range("M32").Value = -561.5
$ws.Range("N32").Value = -2126.75
$ws.Range("H41").Value = 1174.6471
$ws.Range("I41").Value = 1391.8334
$ws.Range("J41").Value = 653.4
$ws.Range("K41").Value = 1391.8334
$ws.Range("L41").Value = 653.4
$ws.Range("M41").Value = -951.8334
$ws.Range("N41").Value = -1533.4
$ws.Range("H55").Value = 380.1
$ws.Range("J55").Value = 400
$ws.Range("L55").Value = 400
$ws.Range("N55").Value = -828
$ws.Range("H98").Value = 1446.7222
$ws.Range("I98").Value = 1516.7858
$ws.Range("K98").Value = 1516.7858
$ws.Range("M98").Value = -18.78580000000011
$ws.Range("H107").Value = 4095.4773
$ws.Range("I107").Value = 5796.069
$ws.Range("J107").Value = 807.6667
$ws.Range("K107").Value = 5796.069
$ws.Range("L107").Value = 807.6667
$ws.Range("M107").Value = -3876.069
$ws.Range("N107").Value = -4647.6667
$ws.Range("H122").Value = 1446.7222
$ws.Range("I122").Value = 1516.7858
$ws.Range("K122").Value = 4550.357400000001
$ws.Range("M122").Value = -2100.357400000001
$ws.Range("H125").Value = 4231.6665
$ws.Range("J125").Value = 1347.5
$ws.Range("L125").Value = 12127.5
$ws.Range("N125").Value = -17047.5
$ws.Range("H135").Value = 17859152
$ws.Range("I135").Value = 1739.091
$ws.Range("J135").Value = 83336340
$ws.Range("K135").Value = 15651.819
$ws.Range("L135").Value = 750027060
$ws.Range("M135").Value = -13116.819
$ws.Range("N135").Value = -750032130

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10018.282
$ws.Range("I32").Value = 12304.697
$ws.Range("K32").Value = 12304.697
$ws.Range("M32").Value = -12017.697
$ws.Range("H110").Value = 1478.4166
$ws.Range("I110").Value = 1249.1818
$ws.Range("J110").Value = 4000
$ws.Range("K110").Value = 1249.1818
$ws.Range("L110").Value = 4000
$ws.Range("M110").Value = 795.8181999999999
$ws.Range("N110").Value = -8090
$ws.Range("H132").Value = 6581166.5
$ws.Range("I132").Value = 7814338
$ws.Range("J132").Value = 4252
$ws.Range("K132").Value = 23443014
$ws.Range("L132").Value = 12756
$ws.Range("M132").Value = -23440484
$ws.Range("N132").Value = -17816

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4335.7427
$ws.Range("I105").Value = 3075.1
$ws.Range("K105").Value = 3075.1
$ws.Range("M105").Value = -1328.1
$ws.Range("H134").Value = 3709.111
$ws.Range("I134").Value = 2197.4285
$ws.Range("J134").Value = 9000
$ws.Range("K134").Value = 6592.2855
$ws.Range("L134").Value = 27000
$ws.Range("M134").Value = -4057.2855
$ws.Range("N134").Value = -32070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7095622
$ws.Range("I31").Value = 3383.432
$ws.Range("J31").Value = 111115120
$ws.Range("K31").Value = 3383.432
$ws.Range("L31").Value = 111115120
$ws.Range("M31").Value = -3088.432
$ws.Range("N31").Value = -111115710
$ws.Range("H34").Value = 7095622
$ws.Range("I34").Value = 3383.432
$ws.Range("J34").Value = 111115120
$ws.Range("K34").Value = 3383.432
$ws.Range("L34").Value = 111115120
$ws.Range("M34").Value = -3181.432
$ws.Range("N34").Value = -111115524
$ws.Range("H62").Value = 2300
$ws.Range("I62").Value = 2300
$ws.Range("K62").Value = 2300
$ws.Range("M62").Value = -1676
$ws.Range("H65").Value = 2300
$ws.Range("I65").Value = 2300
$ws.Range("K65").Value = 11500
$ws.Range("M65").Value = -8380
$ws.Range("H94").Value = 3726.75
$ws.Range("I94").Value = 2480.2
$ws.Range("J94").Value = 4617.143
$ws.Range("K94").Value = 2480.2
$ws.Range("L94").Value = 4617.143
$ws.Range("M94").Value = -2029.2
$ws.Range("N94").Value = -5519.143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4272.727
$ws.Range("J80").Value = 4375
$ws.Range("L80").Value = 13125
$ws.Range("N80").Value = -14997
$ws.Range("H83").Value = 4272.727
$ws.Range("J83").Value = 4375
$ws.Range("L83").Value = 39375
$ws.Range("N83").Value = -48735
$ws.Range("H86").Value = 1288.6666
$ws.Range("I86").Value = 748
$ws.Range("J86").Value = 1559
$ws.Range("K86").Value = 2244
$ws.Range("L86").Value = 4677
$ws.Range("M86").Value = -1058
$ws.Range("N86").Value = -7049
$ws.Range("H89").Value = 1288.6666
$ws.Range("I89").Value = 748
$ws.Range("J89").Value = 1559
$ws.Range("K89").Value = 6732
$ws.Range("L89").Value = 14031
$ws.Range("M89").Value = -804
$ws.Range("N89").Value = -25887
$ws.Range("H121").Value = 1331.2778
$ws.Range("J121").Value = 1583.0714
$ws.Range("L121").Value = 4749.2142
$ws.Range("N121").Value = -7369.2142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3510728.2
$ws.Range("J122").Value = 2249.2222
$ws.Range("L122").Value = 6747.6666
$ws.Range("N122").Value = -11647.6666
$ws.Range("H132").Value = 4092.1052
$ws.Range("I132").Value = 3081.8696
$ws.Range("J132").Value = 5641.1333
$ws.Range("K132").Value = 9245.6088
$ws.Range("L132").Value = 16923.3999
$ws.Range("M132").Value = -6715.6088
$ws.Range("N132").Value = -21983.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6855.722
$ws.Range("I7").Value = 8363
$ws.Range("J7").Value = 5649.9
$ws.Range("K7").Value = 8363
$ws.Range("L7").Value = 5649.9
$ws.Range("M7").Value = -8251
$ws.Range("N7").Value = -5873.9
$ws.Range("H16").Value = 2572.92
$ws.Range("I16").Value = 2101
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 2101
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -1931
$ws.Range("N16").Value = -8340
$ws.Range("H22").Value = 1847.4
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1847.4
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1847.4
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2437.4
$ws.Range("H27").Value = 1847.4
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1847.4
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1847.4
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -2061.4
$ws.Range("H68").Value = 1896.8889
$ws.Range("I68").Value = 1488
$ws.Range("J68").Value = 2224
$ws.Range("K68").Value = 1488
$ws.Range("L68").Value = 2224
$ws.Range("M68").Value = -739
$ws.Range("N68").Value = -3722
$ws.Range("H71").Value = 1896.8889
$ws.Range("I71").Value = 1488
$ws.Range("J71").Value = 2224
$ws.Range("K71").Value = 7440
$ws.Range("L71").Value = 11120
$ws.Range("M71").Value = -3696
$ws.Range("N71").Value = -18608
$ws.Range("H126").Value = 6855.722
$ws.Range("I126").Value = 8363
$ws.Range("J126").Value = 5649.9
$ws.Range("K126").Value = 25089
$ws.Range("L126").Value = 16949.7
$ws.Range("M126").Value = -22619
$ws.Range("N126").Value = -21889.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1952.6562
$ws.Range("I122").Value = 1879.4231
$ws.Range("J122").Value = 2270
$ws.Range("K122").Value = 5638.2693
$ws.Range("L122").Value = 6810
$ws.Range("M122").Value = -3188.2693
$ws.Range("N122").Value = -11710
$ws.Range("H136").Value = 1587.8182
$ws.Range("I136").Value = 1587.8182
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4763.4546
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2213.4546
$ws.Range("N136").ClearContents()
